# Retraining the forecast model for Dragosel Laslea
# Shifts timestamps in column A by +2 days (rows 2-97) and updates
# Notified Production values in column B (rows 2-93) with the
# retrained model's output. Rows 94-97 remain 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Notified Production (MW) values for rows 2-93
$bValues = @{
    2 = 724.54
    3 = 719.773
    4 = 719.586
    5 = 719.097
    6 = 716.026
    7 = 711.239
    8 = 707.9930000000001
    9 = 706.694
    10 = 675.131
    11 = 670.09
    12 = 666.676
    13 = 659.307
    14 = 642.804
    15 = 640.3339999999999
    16 = 647.052
    17 = 651.152
    18 = 646.528
    19 = 644.745
    20 = 644.319
    21 = 647.171
    22 = 667.477
    23 = 672.465
    24 = 680.78
    25 = 681.104
    26 = 698.953
    27 = 702.8630000000001
    28 = 710.482
    29 = 714.982
    30 = 728.373
    31 = 720.95
    32 = 717.526
    33 = 739.006
    34 = 735.3440000000001
    35 = 722.639
    36 = 712.6609999999999
    37 = 699.61
    38 = 631.605
    39 = 630.264
    40 = 630.049
    41 = 628.54
    42 = 633.442
    43 = 639.984
    44 = 644.61
    45 = 650.254
    46 = 685.206
    47 = 691.34
    48 = 697.177
    49 = 701.143
    50 = 723.912
    51 = 727.295
    52 = 731.057
    53 = 733.58
    54 = 755.203
    55 = 756.473
    56 = 758.595
    57 = 762.244
    58 = 818.069
    59 = 829.103
    60 = 838.48
    61 = 850.179
    62 = 975.63
    63 = 1006.756
    64 = 1037.027
    65 = 1068.292
    66 = 1232.764
    67 = 1258.494
    68 = 1286.907
    69 = 1313.416
    70 = 1433.269
    71 = 1453.575
    72 = 1474.513
    73 = 1496.505
    74 = 1590.79
    75 = 1609.126
    76 = 1628.86
    77 = 1649.397
    78 = 1728.262
    79 = 1749.354
    80 = 1767.393
    81 = 1781.468
    82 = 1864.414
    83 = 1871.211
    84 = 1878.781
    85 = 1904.183
    86 = 1937.212
    87 = 1938.842
    88 = 1921.829
    89 = 1925.147
    90 = 1936.867
    91 = 1936.479
    92 = 1938.308
    93 = 1943.87
}

for ($r = 2; $r -le 97; $r++) {
    # Shift the timestamp serial value forward by 2 days
    $oldDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $oldDate + 2

    if ($bValues.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value2 = $bValues[$r]
    }
}
